# "feat: add new task"
#
# Row 45 used to hold a second (duplicate) copy of the
# "zefang-liu/phishing-email-dataset" entry. Replace its contents with a
# brand-new task entry for the "jackhhao/jailbreak-classification" dataset
# (prompt/type columns, new train/test counts, split markers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("B45").Value = "jackhhao/jailbreak-classification"
$ws.Range("F45").Value = "prompt"
$ws.Range("G45").Value = "type"
$ws.Range("H45").Value = 1044
$ws.Range("J45").Value = 262
$ws.Range("K45").Value = "train"
$ws.Range("M45").Value = "test"
